$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2195
$ws1.Range("F3").Value = 914
$ws1.Range("F4").Value = 1661

# Sheet "全部类型" - same rows mirrored with different row offsets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2195
$ws4.Range("F5").Value = 914
$ws4.Range("F6").Value = 1661
